$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update HotStock Top20 table cells (rows 2-21, columns A-C) to the new values
$ws.Range("A2").Value = "华天科技"
$ws.Range("B2").Value = "白银有色"
$ws.Range("C2").Value = "白银有色"
$ws.Range("B3").Value = "常山北明"
$ws.Range("C3").Value = "道生天合"
$ws.Range("A4").Value = "大有能源"
$ws.Range("B4").Value = "闻泰科技"
$ws.Range("C4").Value = "合肥城建"
$ws.Range("A5").Value = "海峡股份"
$ws.Range("B5").Value = "华天科技"
$ws.Range("C5").Value = "华友钴业"
$ws.Range("B6").Value = "三花智控"
$ws.Range("C6").Value = "中兴通讯"
$ws.Range("A7").Value = "海南华铁"
$ws.Range("B7").Value = "合肥城建"
$ws.Range("C7").Value = "山子高科"
$ws.Range("A8").Value = "中兴通讯"
$ws.Range("B8").Value = "山子高科"
$ws.Range("C8").Value = "华天科技"
$ws.Range("A9").Value = "山子高科"
$ws.Range("B9").Value = "大有能源"
$ws.Range("C9").Value = "大有能源"
$ws.Range("A10").Value = "合肥城建"
$ws.Range("B10").Value = "海南华铁"
$ws.Range("C10").Value = "安泰科技"
$ws.Range("A11").Value = "N道生"
$ws.Range("B11").Value = "中兴通讯"
$ws.Range("C11").Value = "常山北明"
$ws.Range("A12").Value = "闻泰科技"
$ws.Range("B12").Value = "N道生"
$ws.Range("C12").Value = "海康威视"
$ws.Range("A13").Value = "三花智控"
$ws.Range("B13").Value = "北方稀土"
$ws.Range("C13").Value = "三花智控"
$ws.Range("A14").Value = "天际股份"
$ws.Range("B14").Value = "海峡股份"
$ws.Range("C14").Value = "海峡股份"
$ws.Range("A15").Value = "东信和平"
$ws.Range("B15").Value = "安泰集团"
$ws.Range("C15").Value = "澄星股份"
$ws.Range("A16").Value = "农业银行"
$ws.Range("B16").Value = "农业银行"
$ws.Range("C16").Value = "闻泰科技"
$ws.Range("B17").Value = "天际股份"
$ws.Range("C17").Value = "天际股份"
$ws.Range("A18").Value = "通富微电"
$ws.Range("B18").Value = "东信和平"
$ws.Range("C18").Value = "东信和平"
$ws.Range("A19").Value = "海峡创新"
$ws.Range("B19").Value = "通富微电"
$ws.Range("C19").Value = "海南华铁"
$ws.Range("A20").Value = "安泰科技"
$ws.Range("B20").Value = "阳光电源"
$ws.Range("C20").Value = "远大控股"
$ws.Range("A21").Value = "深科技"
$ws.Range("B21").Value = "宝泰隆"
$ws.Range("C21").Value = "三孚股份"
